# Update crypto price/volume snapshot (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.247.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "'1.861.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'236.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.4681"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.2867"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "'0.06545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "'21.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.62%  "
$ws.Range("D11").Value = "'0.07928"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "'97.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "'1.867.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'5.183"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "'0.6808"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "'268.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.69%  "
$ws.Range("D17").Value = "'30.235.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "'13.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.97%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "'0.000007410"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").Value = "'2.111.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'5.328"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'6.201"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "'167.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D26").Value = "'9.234"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").Value = "'18.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").Value = "'1.971"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("D29").Value = "'1.384"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("D30").Value = "'0.09926"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("D31").Value = "'4.398"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").Value = "'1.477"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'4.076"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "'0.04703"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'1.136"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.73%  "
$ws.Range("D36").Value = "'0.7051"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "'0.01883"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("E39").Value = "  +4.04%  "
$ws.Range("D40").Value = "'6.252"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("D41").Value = "'74.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("D42").Value = "'1.942"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'0.8484"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'0.4173"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "'0.9996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'103.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "'965.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").Value = "'7.167"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'9.227"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "'34.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "'0.05659"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.34%  "
